$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the two existing data rows (9 and 10) down into rows 11 and 12,
# duplicating their values and formatting.
$ws.Range("A9:FI9").Copy($ws.Range("A11:FI11")) | Out-Null
$ws.Range("A10:FI10").Copy($ws.Range("A12:FI12")) | Out-Null
